# FeSources.xlsx: turn the plain "Link" URLs into clickable-style
# HTML anchor strings, e.g. "<a href='...'>ONS</a>", matching the
# "update data sources to be links" part of the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = "<a href='https://www.ons.gov.uk/peoplepopulationandcommunity/wellbeing/articles/subnationalindicatorsexplorer/2022-01-06'>ONS</a>"
$ws.Range("D4").Value = "<a href='https://www.aoc.co.uk/research-unit/data-sources'>AOC</a>"
$ws.Range("D5").Value = "<a href='https://census.gov.uk/local-authorities'>Census</a>"
$ws.Range("D6").Value = "<a href='https://www.nomisweb.co.uk/'>Nomis</a>"
$ws.Range("D7").Value = "<a href='https://explore-education-statistics.service.gov.uk/'>EES</a>"
